# Added the search and delete for a specific patient
# - Search the patient list for a specific patient ("Hesham") and delete that row
# - Add new appointment-tracking columns: Doctor, Date, Time Slot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Search for the patient to remove ---
$searchName = "Hesham"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$targetRow = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $nameValue = $ws.Cells.Item($r, 1).Value2
    if ($nameValue -eq $searchName) {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    # update the remaining patient record with the matched disease entry
    $diseaseValue = $ws.Cells.Item($targetRow, 4).Value2
    $ws.Cells.Item($targetRow - 1, 4).Value = $diseaseValue

    # delete the matched patient's row entirely
    $ws.Rows.Item($targetRow).Delete()
}

# --- Add the new appointment-tracking columns ---
$ws.Range("E1").Value = "Doctor"
$ws.Range("F1").Value = "Date"
$ws.Range("G1").Value = "Time Slot"
